$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Add two new experiment rows (66, 67) logged the next day, matching the
#    look of rows 25/26 (date in col A formatted like "d-mmm", no fill-in for C:H).
#    These use brand new description text, so enter them first so the new
#    shared strings land before the edited strings below.
$ws.Range("A27").Value = 43374
$ws.Range("A27").NumberFormat = "d-mmm"
$ws.Range("B27").Value = 66
$ws.Range("I27").Value = "Experiment with both GPUs, small complex capsnet, batch size=15, failed"

$ws.Range("A28").Value = 43374
$ws.Range("A28").NumberFormat = "d-mmm"
$ws.Range("B28").Value = 67
$ws.Range("I28").Value = "Experiment with both GPUs, small complex capsnet, batch size=10, failed"

# 2. Append ", batch_size=5" to the existing descriptions for experiments 64 and 65
$ws.Range("I25").Value = "Running on Asia, small complex capsnet, no skips, full data, 10 epochs, batch_size=5"
$ws.Range("I26").Value = "Running on Europe small complex capsnet, with skips, full data, 10 epochs, batch_size=5"

# 3. Move the active selection the way it ended up after the edits
$ws.Range("I30").Select()
